$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B edits - new string order: B6 first, then B3, B4, B5
$ws.Range("B6").Value = "<s><pad></s><unk><mask>"
$ws.Range("B3").Value = "CLS,SEP,MASK, PAD, UNK"
$ws.Range("B4").Value = "CLS,SEP,MASK, PAD, <unk>"
$ws.Range("B5").Value = "CLS,SEP,MASK, PAD, UNK"

# Column G edits - new string order: G2, G3, G4, G5, G6
$ws.Range("G2").Value = "algoritmo"
$ws.Range("G3").Value = "Wordpice"
$ws.Range("G4").Value = "SentencePiece"
$ws.Range("G5").Value = "Wordpice"
$ws.Range("G6").Value = "byte-pair-encoding (BPE) [ Sennrich et al. ] "

# Selection change
$ws.Range("G5").Select()
